$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (req 1.2.2): update requirement text
$ws.Range("B5").Value = 'If username input is empty and login button is clicked display message "Epic sadface: Username is required"'

# Row 9 (req 2.1.1): update requirement text
$ws.Range("B9").Value = "Clicking on add to cart should add item to the cart and display remove button"

# Row 10 (req 2.1.2): text updated (fix "iten" typo + merged with old 2.1.3 wording), and row height increased
$ws.Range("B10").Value = "Clicking on add to cart should add item to the cart and display remove button and clicking remove button should remove item from the cart"
$ws.Rows(10).RowHeight = 25.5

# Row 11 used to be 2.1.3 (now removed) - shift old row 12 (2.2.1) content up into row 11
$ws.Range("A11").Value = "2.2.1"
$ws.Range("B11").Value = 'Clicking on a filter dropdown and selecting "Price (high to low)" should display the most expensive item first'

# Row 12 used to be 2.2.1 - shift old row 13 (2.2.2) content up into row 12
$ws.Range("A12").Value = "2.2.2"
$ws.Range("B12").Value = 'Clicking on a filter dropdown and selecting "Price (low to high)" should display the cheapest item first'

# Row 13 used to be 2.2.2 - now cleared out (empty), reset row height back to default
$ws.Range("A13").Value = ""
$ws.Range("B13").Value = ""
$ws.Rows(13).AutoFit()

# Update the active selection to C10
$ws.Range("C10").Select()
